$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("D3").Value = 44468
$ws.Range("N3").Value = 29000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 29500
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"
$ws.Range("S3").Value = 2950
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44524
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 23000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 23500
$ws.Range("S4").Value = 1958

# Row 5
$ws.Range("D5").Value = 44441
$ws.Range("M5").Value = 100
$ws.Range("Q5").Value = "`$/caja 12 kilos"
$ws.Range("S5").Value = 2458
$ws.Range("T5").Value = 12

# Row 6
$ws.Range("D6").Value = 44160
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 19000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 19500
$ws.Range("Q6").Value = "`$/caja 13 kilos"
$ws.Range("S6").Value = 1500
$ws.Range("T6").Value = 13

# Row 7
$ws.Range("D7").Value = 44496
$ws.Range("L7").Value = "Primera"
$ws.Range("N7").Value = 23000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 23500
$ws.Range("Q7").Value = "`$/caja 12 kilos"
$ws.Range("S7").Value = 1958
$ws.Range("T7").Value = 12

# Row 8
$ws.Range("D8").Value = 44482
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 26000
$ws.Range("P8").Value = 25500
$ws.Range("S8").Value = 2125

# Row 9
$ws.Range("D9").Value = 44475
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 32000
$ws.Range("O9").Value = 33000
$ws.Range("P9").Value = 32500
$ws.Range("S9").Value = 2708

# Row 11
$ws.Range("D11").Value = 44489
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 24500
$ws.Range("S11").Value = 2042
